# Auto-generated edit script: update market-price columns (H-N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 14054.444
$ws.Range("I11").Value = 14054.444
$ws.Range("K11").Value = 14054.444
$ws.Range("M11").Value = -13914.444
$ws.Range("H18").Value = 1363.3636
$ws.Range("I18").Value = 1363.3636
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1363.3636
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1079.3636
$ws.Range("H39").Value = 1219.6666
$ws.Range("I39").Value = 577
$ws.Range("J39").Value = 2023
$ws.Range("K39").Value = 1731
$ws.Range("L39").Value = 6069
$ws.Range("M39").Value = -1435
$ws.Range("N39").Value = -6661
$ws.Range("H40").Value = 1680.375
$ws.Range("J40").Value = 2074
$ws.Range("L40").Value = 2074
$ws.Range("N40").Value = -2424
$ws.Range("H42").Value = 895.0625
$ws.Range("I42").Value = 82.09999999999999
$ws.Range("J42").Value = 2250
$ws.Range("K42").Value = 246.3
$ws.Range("L42").Value = 6750
$ws.Range("M42").Value = -16.29999999999998
$ws.Range("N42").Value = -7210
$ws.Range("H63").Value = 28000
$ws.Range("J63").Value = 28000
$ws.Range("L63").Value = 28000
$ws.Range("N63").Value = -29248
$ws.Range("H66").Value = 28000
$ws.Range("J66").Value = 28000
$ws.Range("L66").Value = 84000
$ws.Range("N66").Value = -90240
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H123").Value = 99498.5
$ws.Range("J123").Value = 99498.5
$ws.Range("L123").Value = 99498.5
$ws.Range("N123").Value = -109298.5
$ws.Range("H132").Value = 999.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 999.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 2998.5
$ws.Range("N132").Value = -8058.5
$ws.Range("H135").Value = 655.3333
$ws.Range("I135").Value = 655.3333
$ws.Range("K135").Value = 5897.9997
$ws.Range("M135").Value = -3362.9997
$ws.Range("H141").Value = 1697.909
$ws.Range("I141").Value = 1540.6666
$ws.Range("K141").Value = 4621.9998
$ws.Range("M141").Value = 558.0002000000004
$ws.Range("N18").Value = $null
$ws.Range("N105").Value = $null
$ws.Range("M132").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10223.223
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10223.223
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10223.223
$ws.Range("N32").Value = -10797.223
$ws.Range("H33").Value = 100000000
$ws.Range("I33").Value = 100000000
$ws.Range("K33").Value = 100000000
$ws.Range("M33").Value = -99999671
$ws.Range("H61").Value = 1053283.1
$ws.Range("I61").Value = 25455.127
$ws.Range("J61").Value = 5444911.5
$ws.Range("K61").Value = 25455.127
$ws.Range("L61").Value = 5444911.5
$ws.Range("M61").Value = -25243.127
$ws.Range("N61").Value = -5445335.5
$ws.Range("H74").Value = 349465.22
$ws.Range("J74").Value = 1013589.25
$ws.Range("L74").Value = 1013589.25
$ws.Range("N74").Value = -1015337.25
$ws.Range("H77").Value = 349465.22
$ws.Range("J77").Value = 1013589.25
$ws.Range("L77").Value = 5067946.25
$ws.Range("N77").Value = -5076682.25
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H132").Value = 1480.5902
$ws.Range("I132").Value = 1001.39624
$ws.Range("K132").Value = 3004.18872
$ws.Range("M132").Value = -474.1887200000001
$ws.Range("H136").Value = 1053283.1
$ws.Range("I136").Value = 25455.127
$ws.Range("J136").Value = 5444911.5
$ws.Range("K136").Value = 76365.38099999999
$ws.Range("L136").Value = 16334734.5
$ws.Range("M136").Value = -73815.38099999999
$ws.Range("N136").Value = -16339834.5
$ws.Range("M32").Value = $null
$ws.Range("N92").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 24000
$ws.Range("J62").Value = 24000
$ws.Range("L62").Value = 24000
$ws.Range("N62").Value = -25372
$ws.Range("H65").Value = 24000
$ws.Range("J65").Value = 24000
$ws.Range("L65").Value = 72000
$ws.Range("N65").Value = -78864
$ws.Range("H105").Value = 7316.515
$ws.Range("I105").Value = 6393.609
$ws.Range("K105").Value = 6393.609
$ws.Range("M105").Value = -4646.609
$ws.Range("H134").Value = 21430138
$ws.Range("I134").Value = 1236.1852
$ws.Range("K134").Value = 3708.5556
$ws.Range("M134").Value = -1173.5556

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 856.8570999999999
$ws.Range("I22").Value = 946.5
$ws.Range("J22").Value = 737.3333
$ws.Range("K22").Value = 946.5
$ws.Range("L22").Value = 737.3333
$ws.Range("M22").Value = -596.5
$ws.Range("N22").Value = -1437.3333
$ws.Range("H33").Value = 6803.7896
$ws.Range("I33").Value = 1300
$ws.Range("J33").Value = 12919.111
$ws.Range("K33").Value = 1300
$ws.Range("L33").Value = 12919.111
$ws.Range("M33").Value = -921
$ws.Range("N33").Value = -13677.111
$ws.Range("H75").Value = 77500
$ws.Range("J75").Value = 77500
$ws.Range("L75").Value = 77500
$ws.Range("N75").Value = -79496
$ws.Range("H78").Value = 77500
$ws.Range("J78").Value = 77500
$ws.Range("L78").Value = 232500
$ws.Range("N78").Value = -242484
$ws.Range("H134").Value = 2754.4
$ws.Range("I134").Value = 2499.1538
$ws.Range("J134").Value = 3228.4285
$ws.Range("K134").Value = 7497.4614
$ws.Range("L134").Value = 9685.2855
$ws.Range("M134").Value = -4962.4614
$ws.Range("N134").Value = -14755.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1206.6666
$ws.Range("I8").Value = 1206.6666
$ws.Range("K8").Value = 3619.9998
$ws.Range("M8").Value = -3480.9998
$ws.Range("H12").Value = 438
$ws.Range("I12").Value = 196.2
$ws.Range("J12").Value = 531
$ws.Range("K12").Value = 588.5999999999999
$ws.Range("L12").Value = 1593
$ws.Range("M12").Value = -415.5999999999999
$ws.Range("N12").Value = -1939
$ws.Range("H17").Value = 964
$ws.Range("I17").Value = 510
$ws.Range("J17").Value = 1266.6666
$ws.Range("K17").Value = 1530
$ws.Range("L17").Value = 3799.9998
$ws.Range("M17").Value = -1361
$ws.Range("N17").Value = -4137.9998
$ws.Range("H111").Value = 1700
$ws.Range("J111").Value = 1400
$ws.Range("L111").Value = 4200
$ws.Range("N111").Value = -10334
$ws.Range("H119").Value = 71440780
$ws.Range("I119").Value = 100007500
$ws.Range("K119").Value = 300022500
$ws.Range("M119").Value = -300017662
$ws.Range("H122").Value = 9526329
$ws.Range("I122").Value = 16667103
$ws.Range("J122").Value = 5297
$ws.Range("K122").Value = 150003927
$ws.Range("L122").Value = 47673
$ws.Range("M122").Value = -150001477
$ws.Range("N122").Value = -52573
$ws.Range("H131").Value = 4682.3335
$ws.Range("J131").Value = 6328.3335
$ws.Range("L131").Value = 18985.0005
$ws.Range("N131").Value = -29065.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4113.375
$ws.Range("I122").Value = 4063.1538
$ws.Range("K122").Value = 12189.4614
$ws.Range("M122").Value = -9739.4614
$ws.Range("H132").Value = 1172055.9
$ws.Range("I132").Value = 14663
$ws.Range("J132").Value = 1833423.1
$ws.Range("K132").Value = 43989
$ws.Range("L132").Value = 5500269.300000001
$ws.Range("M132").Value = -41459
$ws.Range("N132").Value = -5505329.300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 4004
$ws.Range("I17").Value = 4004
$ws.Range("K17").Value = 4004
$ws.Range("M17").Value = -3834
$ws.Range("H122").Value = 3495.9355
$ws.Range("I122").Value = 3128.36
$ws.Range("K122").Value = 9385.08
$ws.Range("M122").Value = -6935.08
$ws.Range("H136").Value = 55206.42
$ws.Range("I136").Value = 93504
$ws.Range("K136").Value = 280512
$ws.Range("M136").Value = -277962

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 28000.6
$ws.Range("J4").Value = 13334.333
$ws.Range("L4").Value = 13334.333
$ws.Range("N4").Value = -13560.333
$ws.Range("H6").Value = 3749.75
$ws.Range("I6").Value = 4666.3335
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 4666.3335
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = -4551.3335
$ws.Range("N6").Value = -1230
$ws.Range("I107").Value = 1999.5
$ws.Range("J107").Value = 5717066.5
$ws.Range("K107").Value = 5998.5
$ws.Range("L107").Value = 17151199.5
$ws.Range("M107").Value = -4078.5
$ws.Range("N107").Value = -17155039.5
$ws.Range("H113").Value = 649.7778
$ws.Range("I113").Value = 577
$ws.Range("K113").Value = 1731
$ws.Range("M113").Value = 439
